# Apply latest cryptos price / volume(1h) snapshot (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.451.16"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "'2.590.63"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'507.03"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'153.89"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -6.60%  "
$ws.Range("D9").Value = "'2.599.93"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'6.64"
$ws.Range("E10").Value = "  +8.01%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'0.346"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "'0.128"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'3.045.79"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'60.447.95"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "'21.64"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "'2.591.91"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'4.82"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "'347.34"
$ws.Range("E20").Value = "  +3.62%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'60.29"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'0.421"
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "'2.700.98"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "'0.0₃0847"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'19.34"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").Value = "'154.47"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").Value = "'5.76"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").Value = "'4.01"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  +18.41%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "'3.77"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").Value = "'35.86"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").Value = "'297.26"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "'0.619"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").Value = "'0.0999"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'19.67"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").Value = "'4.93"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.29"
$ws.Range("E51").Value = "  +0.14%  "
